$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Fill in row 8 with the new "features_custom_2" result row
$ws.Range("A8").Value = "features_custom_2"
$ws.Range("B8").Value = "76,0.573"
$ws.Range("C8").Value = "102,0.4269"
$ws.Range("D8").Value = "103,0.4213"
$ws.Range("E8").Value = "90,0.4943"

# Update the active selection on the sheet to E11
$ws.Range("E11").Select()
